# "added textbox visible flag"
# The diff removes the "TextBox 1" shape from every slide's shape tree and,
# as a consequence of the picture/shape generator re-emitting the remaining
# shapes, every shape that came after it is renumbered/renamed down by one
# (Picture 2 -> Picture 1, Picture 3 -> Picture 2, ..., Rounded Rectangle 10
# -> Rounded Rectangle 9). Reproduce both effects here: delete the textbox,
# then rename the remaining shapes to match.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    # Find and remove the "TextBox 1" shape (always id=2 / first shape in
    # these decks), wherever it happens to sit in the shape collection.
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Name -eq "TextBox 1") {
            $shp.Delete()
        }
    }

    # Shift every "Picture N" / "Rounded Rectangle N" shape name down by one
    # to fill the numbering gap left by the removed textbox.
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        $name = $shp.Name

        if ($name -match "^(Picture|Rounded Rectangle) (\d+)$") {
            $prefix = $matches[1]
            $num = [int]$matches[2]
            $shp.Name = "$prefix $($num - 1)"
        }
    }
}
